# Add two new "Scene" configure fields (CamOffestPos, CamOffestRot), each
# laid out exactly like the existing fields above them (e.g. row 12 /
# "LoadingUI"): A=Id, B=Type, C..F=bool flags, G/H=numbers,
# I=RelationValue, J=Desc.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$newFields = @(
    @{ Row = 13; Id = "CamOffestPos" },
    @{ Row = 14; Id = "CamOffestRot" }
)

foreach ($field in $newFields) {
    $r = $field.Row

    $ws.Cells.Item($r, 1).Value = $field.Id    # A: Id
    $ws.Cells.Item($r, 2).Value = "string"     # B: Type
    $ws.Cells.Item($r, 3).Value = $false       # C: Public
    $ws.Cells.Item($r, 4).Value = $false       # D: Private
    $ws.Cells.Item($r, 5).Value = $false       # E: Save
    $ws.Cells.Item($r, 6).Value = $true        # F: View
    $ws.Cells.Item($r, 7).Value = 0            # G: Index
    $ws.Cells.Item($r, 8).Value = 0            # H: SaveInterval
    $ws.Cells.Item($r, 9).Value = "Friend"     # I: RelationValue
    $ws.Cells.Item($r, 10).Value = "acctorid"  # J: Desc

    # Columns A, B, I, J carry the sheet's "text" cell style (same as every
    # other data row); re-applying the text number format makes the new
    # cells reuse that existing style instead of minting a new one.
    foreach ($col in @(1, 2, 9, 10)) {
        $ws.Cells.Item($r, $col).NumberFormat = "@"
    }
}

$ws.Range("A14").Select()
